# semana 20 de 2024
# Update the "Esperado" (C), "Observado" (D) and "valor p" (E) columns
# in poisson.xlsx with the refreshed weekly figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - 113 Desnutricion aguda en menores de 5 anos
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = 0.02

# Row 4 - 115 Cancer en menores de 18 anos
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 1

# Row 5 - 155 Cancer de la mama y cuello uterino
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 11
$ws.Range("E5").Value = 0.01

# Row 6 - 210 Dengue
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 104

# Row 9 - 220 Dengue grave
$ws.Range("D9").Value = 2

# Row 11 - 300 Agresiones por animales potencialmente transmisores de rabia
$ws.Range("C11").Value = 42
$ws.Range("D11").Value = 31
$ws.Range("E11").Value = 0.01

# Row 14 - 342 Enfermedades huerfanas - raras
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = 0.22

# Row 15 - 346 Ira por virus nuevo
$ws.Range("C15").Value = 7

# Row 16 - 348 Infeccion respiratoria aguda grave irag inusitada
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 0.37

# Row 18 - 355 Enfermedad transmitida por alimentos o agua (eta)
$ws.Range("C18").Value = 1
$ws.Range("E18").Value = 0.37

# Row 20 - 357 Iad - infecciones asociadas a dispositivos - individual
$ws.Range("C20").Value = 2
$ws.Range("E20").Value = 0.14

# Row 21 - 365 Intoxicaciones
$ws.Range("C21").Value = 6
$ws.Range("D21").Value = 3

# Row 24 - 455 Leptospirosis
$ws.Range("D24").Value = 1

# Row 25 - 465 Malaria
$ws.Range("D25").Value = 1

# Row 27 - 549 Morbilidad materna extrema
$ws.Range("C27").Value = 7
$ws.Range("D27").Value = 14
$ws.Range("E27").Value = 0.01

# Row 28 - 560 Mortalidad perinatal y neonatal tardia
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0.14

# Row 29 - 580 Mortalidad por dengue
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = 0

# Row 30 - 620 Parotiditis
$ws.Range("D30").Value = 6
$ws.Range("E30").Value = 0

# Row 31 - 740 Sifilis congenita
$ws.Range("D31").Value = 1
$ws.Range("E31").Value = 0

# Row 32 - 750 Sifilis gestacional
$ws.Range("D32").Value = 1

# Row 33 - 813 Tuberculosis
$ws.Range("C33").Value = 9
$ws.Range("D33").Value = 6
$ws.Range("E33").Value = 0.09

# Row 34 - 831 Varicela individual
$ws.Range("C34").Value = 9
$ws.Range("D34").Value = 4
$ws.Range("E34").Value = 0.03

# Row 35 - 850 Vih/sida/mortalidad por sida
$ws.Range("C35").Value = 9
$ws.Range("E35").Value = 0.03
